# Apply the Univolei "rallies" sheet update:
#  - "sets" sheet: set_id 3 (row 4) home_points (D4) goes from 17 to 22
#  - "rallies" sheet: append 5 new rally rows (81-85) that close out set 3
#    (home team "NOS" scoring points 18-22)

$wb = $excel.ActiveWorkbook

# --- 1. Update the "sets" sheet: home_points for the 3rd set is now 22 ---
$sets = $wb.Worksheets.Item("sets")
$sets.Range("D4").Value = 22

# --- 2. Append the new rally rows to the "rallies" sheet ---
$rallies = $wb.Worksheets.Item("rallies")

# Each tuple: rally_id, match_id, set_number, rally_no, side, position,
#             player_number, action, result, who_scored, score_home,
#             score_away, raw_text, position_zone, pos_fb, frente_fundo
$newRows = @(
    @(80, 1, 3, 18, "NOS", 5, "LOB",     "PONTO", "NOS", 18, 0, "1 5 lob", "FRENTE", "FRENTE", "FRENTE"),
    @(81, 1, 3, 19, "NOS", 5, "SEGUNDA", "PONTO", "NOS", 19, 0, "1 5 seg", "FRENTE", "FRENTE", "FRENTE"),
    @(82, 1, 3, 20, "NOS", 5, "SEGUNDA", "PONTO", "NOS", 20, 0, "1 5 seg", "FRENTE", "FRENTE", "FRENTE"),
    @(83, 1, 3, 21, "NOS", 6, "PIPE",    "PONTO", "NOS", 21, 0, "1 6 pi",  "FRENTE", "FRENTE", "FRENTE"),
    @(84, 1, 3, 22, "NOS", 5, "LOB",     "PONTO", "NOS", 22, 0, "1 5 lob", "FRENTE", "FRENTE", "FRENTE")
)

$startRow = 81
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $data = $newRows[$i]

    $rallies.Cells.Item($r, 1).Value = $data[0]    # A rally_id
    $rallies.Cells.Item($r, 2).Value = $data[1]    # B match_id
    $rallies.Cells.Item($r, 3).Value = $data[2]    # C set_number
    $rallies.Cells.Item($r, 4).Value = $data[3]    # D rally_no
    $rallies.Cells.Item($r, 5).Value = $data[4]    # E side
    # F (position) is left blank for these rows
    $rallies.Cells.Item($r, 7).Value = $data[5]    # G player_number
    $rallies.Cells.Item($r, 8).Value = $data[6]    # H action
    $rallies.Cells.Item($r, 9).Value = $data[7]    # I result
    $rallies.Cells.Item($r, 10).Value = $data[8]   # J who_scored
    $rallies.Cells.Item($r, 11).Value = $data[9]   # K score_home
    $rallies.Cells.Item($r, 12).Value = $data[10]  # L score_away
    $rallies.Cells.Item($r, 13).Value = $data[11]  # M raw_text
    $rallies.Cells.Item($r, 14).Value = $data[12]  # N position_zone
    $rallies.Cells.Item($r, 15).Value = $data[13]  # O pos_fb
    $rallies.Cells.Item($r, 16).Value = $data[14]  # P frente_fundo
}
